# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracker.
#
# The sheet keeps a running log where every row's Day cell (column A) is
# formatted as a full date-time EXCEPT the most recent row, which is
# formatted as a plain date. When a new day's data is appended, the
# previous "last row" gets switched to the date-time format and the new
# row takes over the plain-date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the current last row of data (the row right below the header).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# The previous last row's date switches from date-only to date+time format,
# matching every other historical row.
$ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row's data.
$newRow = $lastRow + 1
$prevDay = $ws.Cells.Item($lastRow, 1).Value2

$ws.Cells.Item($newRow, 1).Value = $prevDay + 1
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($newRow, 2).Value = 121
$ws.Cells.Item($newRow, 3).Value = 121
$ws.Cells.Item($newRow, 4).Value = 124
